# Hortaliza, Vega Modelo de Temuco - Berenjena: weekly data refresh.
# A new observation row is inserted right before the current row 263,
# pushing the existing rows 263-323 down to 264-324 (dimension grows to
# A1:R324). All other rows are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 263 (shifts rows 263..323 -> 264..324).
$ws.Rows.Item(263).Insert()

# Populate the newly inserted row with the new data point.
$ws.Cells.Item(263, 1).Value  = 10
$ws.Cells.Item(263, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(263, 3).Value  = "La Araucanía"
$ws.Cells.Item(263, 4).Value  = 44782
$ws.Cells.Item(263, 5).Value  = 9
$ws.Cells.Item(263, 6).Value  = 100112001
$ws.Cells.Item(263, 7).Value  = "Berenjena"
$ws.Cells.Item(263, 8).Value  = "Sin especificar"
$ws.Cells.Item(263, 9).Value  = "Primera"
$ws.Cells.Item(263, 10).Value = 50
$ws.Cells.Item(263, 11).Value = 14000
$ws.Cells.Item(263, 12).Value = 15000
$ws.Cells.Item(263, 13).Value = 14600
$ws.Cells.Item(263, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(263, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(263, 16).Value = 243
$ws.Cells.Item(263, 17).Value = 60
$ws.Cells.Item(263, 18).Value = "Hortaliza"
